$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Preserve the D:E columns as plain text (they hold numeric-looking strings
# like "243.59" and percentages like "-0.18%") so Excel does not auto-convert
# them into numeric/percentage values.
$ws.Range("D2:E50").NumberFormat = "@"

$ws.Range("D2").Value = "243.59"
$ws.Range("E2").Value = "-0.18%"
$ws.Range("D3").Value = "29.82"
$ws.Range("E3").Value = "13.24%"
$ws.Range("D4").Value = "5.156"
$ws.Range("E4").Value = "0.56%"
$ws.Range("D5").Value = "0.05662"
$ws.Range("E5").Value = "1.22%"
$ws.Range("D6").Value = "6.524"
$ws.Range("E6").Value = "0.78%"
$ws.Range("D7").Value = "0.8455"
$ws.Range("E7").Value = "2.85%"
$ws.Range("D8").Value = "0.8636"
$ws.Range("E8").Value = "2.98%"
$ws.Range("B9").Value = "WazirX"
$ws.Range("C9").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D9").Value = "0.1337"
$ws.Range("E9").Value = "0.68%"
$ws.Range("B10").Value = "MandalaExchangeToken"
$ws.Range("C10").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D10").Value = "0.06907"
$ws.Range("E10").Value = "-1.26%"
$ws.Range("B11").Value = "BitrueCoin"
$ws.Range("C11").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D11").Value = "0.02908"
$ws.Range("E11").Value = "0.82%"
$ws.Range("B12").Value = "BitMartToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D12").Value = "0.09382"
$ws.Range("E12").Value = "-0.03%"
$ws.Range("B13").Value = "BitForexToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D13").Value = "0.001521"
$ws.Range("E13").Value = "-0.53%"
$ws.Range("B14").Value = "CoinExToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D14").Value = "0.04166"
$ws.Range("E14").Value = "-10.09%"
$ws.Range("B15").Value = "TigerCash"
$ws.Range("C15").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D15").Value = "0.006200"
$ws.Range("E15").Value = "-0.35%"
$ws.Range("B16").Value = "LEO"
$ws.Range("C16").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D16").Value = "3.507"
$ws.Range("E16").Value = "-4.02%"
$ws.Range("B17").Value = "GateToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D17").Value = "3.022"
$ws.Range("E17").Value = "-0.31%"
$ws.Range("B18").Value = "BTSEToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D18").Value = "2.128"
$ws.Range("E18").Value = "-2.50%"
$ws.Range("B19").Value = "BitpandaEcosystemToken"
$ws.Range("C19").Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
$ws.Range("D19").Value = "0.3149"
$ws.Range("E19").Value = "1.20%"
$ws.Range("B20").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C20").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D20").Value = "0.03256"
$ws.Range("E20").Value = "4.05%"
$ws.Range("B21").Value = "ProBitToken"
$ws.Range("C21").Value = "https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
$ws.Range("D21").Value = "0.1303"
$ws.Range("E21").Value = "0.29%"
$ws.Range("B22").Value = "MCDex"
$ws.Range("C22").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D22").Value = "3.602"
$ws.Range("E22").Value = "-3.58%"
$ws.Range("B23").Value = "ZBToken"
$ws.Range("C23").Value = "https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb"
$ws.Range("D23").Value = "0.1373"
$ws.Range("E23").Value = "-0.06%"
$ws.Range("B24").Value = "One"
$ws.Range("C24").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D24").Value = "0.0005977"
$ws.Range("E24").Value = "-0.20%"
$ws.Range("D25").Value = "0.001209"
$ws.Range("E25").Value = "-2.70%"
$ws.Range("E26").Value = "-1.26%"
$ws.Range("D27").Value = "0.0001180"
$ws.Range("E28").Value = "0.30%"
$ws.Range("E40").Value = "1.90%"
$ws.Range("D41").Value = "0.005328"
$ws.Range("E41").Value = "-13.23%"
$ws.Range("D42").Value = "0.1058"
$ws.Range("E42").Value = "0.75%"
$ws.Range("E43").Value = "-3.75%"
$ws.Range("D44").Value = "0.009783"
$ws.Range("E44").Value = "20.14%"
$ws.Range("D45").Value = "0.00005092"
$ws.Range("E45").Value = "-4.74%"
$ws.Range("E46").Value = "-0.05%"
$ws.Range("D47").Value = "0.09995"
$ws.Range("E47").Value = "-30.57%"
$ws.Range("D48").Value = "0.002819"
$ws.Range("E48").Value = "22.61%"
$ws.Range("E49").Value = "-0.05%"
$ws.Range("E50").Value = "-0.05%"

# Restore the default (un-styled) cell style now that the text values are set,
# so no stray style index is left behind on these cells.
$ws.Range("D2:E50").Style = "Normal"

